# -----------------------------------------------------------------------
# Edit: fix: Remove unconfirmed expenses, keep only guaranteed figures
#
# - Costs_Tracker sheet:
#     * Add VAT call-outs / "No VAT" / "VAT exempt" notes to several rows
#     * Correct rounded Budgeted/Forecast totals to their precise VAT-
#       inclusive figures
#     * Remove the "Warehouse - Racking & Setup" and "Working Capital"
#       placeholder rows
# - Monthly_Cashflow sheet:
#     * Remove the unconfirmed "Racking & Setup" and "Warehouse Equipment"
#       January entries
#     * Rename "Business Rates (Dec)" to "Business Rates"
#     * Recalculate the Running Balance column for all rows following the
#       removed entries
# -----------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# ======================================================================
# Sheet: Costs_Tracker
# ======================================================================
$costs = $wb.Worksheets.Item("Costs_Tracker")

# --- Row 2: Warehouse - Rent Deposit (7 months) ---
$costs.Range("B2").Value = 154065.6
$costs.Range("D2").Value = 154065.6
$costs.Range("E2").Value = "£128,388 + VAT (£25,677.60) - Due Dec 19th"

# --- Row 3: Warehouse - Q1 Rent ---
$costs.Range("B3").Value = 38516.4
$costs.Range("D3").Value = 38516.4
$costs.Range("E3").Value = "£32,097 + VAT (£6,419.40) - Due Dec 19th"

# --- Row 4: Warehouse - Service Charge (Quarterly) ---
$costs.Range("E4").Value = "£12,000 per quarter - No VAT"

# --- Row 5: Warehouse - Insurance (Annual) ---
$costs.Range("E5").Value = "Yearly upfront - VAT exempt - Due Dec 19th"

# --- Row 6: Warehouse - Business Rates (Monthly) ---
$costs.Range("E6").Value = "£5,000 per month - No VAT"

# --- Row 7: Warehouse - Legal/Professional ---
$costs.Range("B7").Value = 13116.6
$costs.Range("D7").Value = 13116.6
$costs.Range("E7").Value = "£10,930.50 + VAT (£2,186.10) - Due Dec 19th"

# --- Remove unconfirmed rows ---
# Row 10: "Working Capital" placeholder - delete first so row indices
# above it are unaffected while we still need them.
$costs.Rows.Item(10).Delete()
# Row 8: "Warehouse - Racking & Setup" - after deleting this row, the
# former row 9 ("Showroom - Total Completion") shifts up to become row 8.
$costs.Rows.Item(8).Delete()

# ======================================================================
# Sheet: Monthly_Cashflow
# ======================================================================
$cash = $wb.Worksheets.Item("Monthly_Cashflow")

# Remove the unconfirmed January entries: "Warehouse Equipment" (row 11)
# then "Racking & Setup" (row 10). Delete bottom-up so row numbers stay
# valid between calls.
$cash.Rows.Item(11).Delete()
$cash.Rows.Item(10).Delete()

# Rename Dec 2025 business rates line item
$cash.Range("B8").Value = "Business Rates"

# Recalculate the Running Balance column from row 8 downward now that the
# two removed entries (£25,000 + £5,000 gross) no longer reduce it.
$cash.Range("G8").Value = 371824
$cash.Range("G9").Value = 366824
$cash.Range("G10").Value = 361824
$cash.Range("G11").Value = 356824
$cash.Range("G12").Value = 318307.6
$cash.Range("G13").Value = 306307.6
$cash.Range("G14").Value = 301307.6
$cash.Range("G15").Value = 296307.6
$cash.Range("G16").Value = 291307.6
$cash.Range("G17").Value = 252791.2
$cash.Range("G18").Value = 240791.2
